$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Append two new data rows (18 and 19) for KNN (k=1), mirroring the
#     existing per-algorithm row-pair pattern already used for rows 2:17 ---

# Row 18
$ws.Cells.Item(18, 1).Value = "KNN"
$ws.Cells.Item(18, 2).Value = 1
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 5).Value = 2810
$ws.Cells.Item(18, 6).Value = 2755
$ws.Cells.Item(18, 7).Formula = "=E18-F18"
$ws.Cells.Item(18, 8).Formula = "=ROUND(F18*100/E18,4)"
$ws.Cells.Item(18, 9).Formula = "=ROUND(AVERAGE(H18,H19),4)"

# Row 19
$ws.Cells.Item(19, 1).Value = "KNN"
$ws.Cells.Item(19, 2).Value = 1
$ws.Cells.Item(19, 3).Value = 2
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 5).Value = 2810
$ws.Cells.Item(19, 6).Value = 2767
$ws.Cells.Item(19, 7).Formula = "=E19-F19"
$ws.Cells.Item(19, 8).Formula = "=ROUND(F19*100/E19,4)"

# Copy formatting from the row above (17) so the new rows match the
# existing look (centred "Algorithm" column, centred merged "avg" column)
$ws.Range("A17").Copy()
$ws.Range("A18:A19").PasteSpecial(-4122)

$ws.Range("I17").Copy()
$ws.Range("I18:I19").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Merge I18:I19, matching the merge pattern used for every other row pair
$ws.Range("I18:I19").Merge()

# --- Update dimension / selection to mirror scrolling to the newly added pair ---
$ws.Range("I18:I19").Select()

$wb.Save()
